$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "(SITE)Página: LOGIN" task (row 7) is now finished:
#   STATUS goes from "Em processo" to "FEITO" and % CONCLUÍDO goes from 0 to 100%.
$ws.Range("C7").Value = "FEITO"
$ws.Range("F7").Value = 1

# Leave the cursor where the author left it when they saved the file.
$ws.Range("F18").Select()
